$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "List1" (sheet1.xml): results by category, columns B (kategorija), C (grupa), D (staza), E (kontrole) ---
$ws1.Range("B2").Value = 'M12'
$ws1.Range("C2").Value = 'M12; Ž12; Open Short'
$ws1.Range("D2").Value = '1.66–50–9'
$ws1.Range("E2").Value = '31–32–33–34–35–36–37–38–100––––––––––––––––––––'
$ws1.Range("B3").Value = 'OPEN kratka'
$ws1.Range("C3").Value = 'M12; Ž12; Open Short'
$ws1.Range("D3").Value = '1.66–50–9'
$ws1.Range("E3").Value = '31–32–33–34–35–36–37–38–100––––––––––––––––––––'
$ws1.Range("B4").Value = 'Ž12'
$ws1.Range("C4").Value = 'M12; Ž12; Open Short'
$ws1.Range("D4").Value = '1.66–50–9'
$ws1.Range("E4").Value = '31–32–33–34–35–36–37–38–100––––––––––––––––––––'
$ws1.Range("B5").Value = 'M21A'
$ws1.Range("C5").Value = 'M21A'
$ws1.Range("D5").Value = '4.89–220–18'
$ws1.Range("E5").Value = '40–41–42–43–44–45–46–47–48–49–50–51–39–52–53–44–54–100–––––––––––'
$ws1.Range("B6").Value = 'M20'
$ws1.Range("C6").Value = 'M20; M35; Ž21A'
$ws1.Range("D6").Value = '3.7–150–16'
$ws1.Range("E6").Value = '40–41–42–43–58–44–45–46–47–48–61–39–62–60–67–100–––––––––––––'
$ws1.Range("B7").Value = 'M35'
$ws1.Range("C7").Value = 'M20; M35; Ž21A'
$ws1.Range("D7").Value = '3.7–150–16'
$ws1.Range("E7").Value = '40–41–42–43–58–44–45–46–47–48–61–39–62–60–67–100–––––––––––––'
$ws1.Range("B8").Value = 'Ž21A'
$ws1.Range("C8").Value = 'M20; M35; Ž21A'
$ws1.Range("D8").Value = '3.7–150–16'
$ws1.Range("E8").Value = '40–41–42–43–58–44–45–46–47–48–61–39–62–60–67–100–––––––––––––'
$ws1.Range("B9").Value = 'Ž16'
$ws1.Range("C9").Value = 'Ž16; Ž21B'
$ws1.Range("D9").Value = '2.11–80–10'
$ws1.Range("E9").Value = '40–53–44–45–59–51–39–52–60–100–––––––––––––––––––'
$ws1.Range("B10").Value = 'Ž21B'
$ws1.Range("C10").Value = 'Ž16; Ž21B'
$ws1.Range("D10").Value = '2.11–80–10'
$ws1.Range("E10").Value = '40–53–44–45–59–51–39–52–60–100–––––––––––––––––––'
$ws1.Range("B11").Value = 'M14'
$ws1.Range("C11").Value = 'M14; Ž14'
$ws1.Range("D11").Value = '1.83–55–9'
$ws1.Range("E11").Value = '55–40–33–53–35–56–37–67–100––––––––––––––––––––'
$ws1.Range("B12").Value = 'Ž14'
$ws1.Range("C12").Value = 'M14; Ž14'
$ws1.Range("D12").Value = '1.83–55–9'
$ws1.Range("E12").Value = '55–40–33–53–35–56–37–67–100––––––––––––––––––––'
$ws1.Range("B13").Value = 'M55'
$ws1.Range("C13").Value = 'M55; Ž45'
$ws1.Range("D13").Value = '2.54–90–13'
$ws1.Range("E13").Value = '57–40–53–59–51–46–61–47–39–62–60–67–100––––––––––––––––'
$ws1.Range("B14").Value = 'Ž45'
$ws1.Range("C14").Value = 'M55; Ž45'
$ws1.Range("D14").Value = '2.54–90–13'
$ws1.Range("E14").Value = '57–40–53–59–51–46–61–47–39–62–60–67–100––––––––––––––––'
$ws1.Range("B15").Value = 'M45'
$ws1.Range("C15").Value = 'M45; Ž20; Ž35'
$ws1.Range("D15").Value = '2.82–110–11'
$ws1.Range("E15").Value = '57–40–63–62–46–61–64–47–51–45–100––––––––––––––––––'
$ws1.Range("B16").Value = 'Ž20'
$ws1.Range("C16").Value = 'M45; Ž20; Ž35'
$ws1.Range("D16").Value = '2.82–110–11'
$ws1.Range("E16").Value = '57–40–63–62–46–61–64–47–51–45–100––––––––––––––––––'
$ws1.Range("B17").Value = 'Ž35'
$ws1.Range("C17").Value = 'M45; Ž20; Ž35'
$ws1.Range("D17").Value = '2.82–110–11'
$ws1.Range("E17").Value = '57–40–63–62–46–61–64–47–51–45–100––––––––––––––––––'
$ws1.Range("B18").Value = 'M16'
$ws1.Range("C18").Value = 'M16; M21B; Open Long'
$ws1.Range("D18").Value = '2.24–90–11'
$ws1.Range("E18").Value = '57–43–58–44–53–59–51–39–52–60–100––––––––––––––––––'
$ws1.Range("B19").Value = 'M21B'
$ws1.Range("C19").Value = 'M16; M21B; Open Long'
$ws1.Range("D19").Value = '2.24–90–11'
$ws1.Range("E19").Value = '57–43–58–44–53–59–51–39–52–60–100––––––––––––––––––'
$ws1.Range("B20").Value = 'OPEN duga'
$ws1.Range("C20").Value = 'M16; M21B; Open Long'
$ws1.Range("D20").Value = '2.24–90–11'
$ws1.Range("E20").Value = '57–43–58–44–53–59–51–39–52–60–100––––––––––––––––––'
$ws1.Range("B21").Value = 'M16'
$ws1.Range("C21").Value = 'M16; M21B; Open Long'
$ws1.Range("D21").Value = '2.24–90–11'
$ws1.Range("E21").Value = '57–43–58–44––59–51–39–52–60–100––––––––––––––––––'
$ws1.Range("B22").Value = 'Ž65'
$ws1.Range("C22").Value = 'Ž65; Ž70'
$ws1.Range("D22").Value = '1.75–80–9'
$ws1.Range("E22").Value = '65–66–44–45–51–62–63–67–100––––––––––––––––––––'
$ws1.Range("B23").Value = 'M65'
$ws1.Range("C23").Value = 'M65; M70; Ž55'
$ws1.Range("D23").Value = '2.03–80–10'
$ws1.Range("E23").Value = '65–66–53–59–51–39–52–63–67–100–––––––––––––––––––'
$ws1.Range("B24").Value = 'M70'
$ws1.Range("C24").Value = 'M65; M70; Ž55'
$ws1.Range("D24").Value = '2.03–80–10'
$ws1.Range("E24").Value = '65–66–53–59–51–39–52–63–67–100–––––––––––––––––––'
$ws1.Range("B25").Value = 'Ž55'
$ws1.Range("C25").Value = 'M65; M70; Ž55'
$ws1.Range("D25").Value = '2.03–80–10'
$ws1.Range("E25").Value = '65–66–53–59–51–39–52–63–67–100–––––––––––––––––––'

# --- Sheet "List2" (sheet2.xml): ranking summary, columns B/C/E (etapa) and G/H (ukupno) ---
$ws2.Range("B2").Value = 'M21A'
$ws2.Range("C2").Value = 'M21A'
$ws2.Range("E2").Value = 1
$ws2.Range("G2").Value = 'M21A'
$ws2.Range("H2").Value = 1
$ws2.Range("B3").Value = 'Ž21A'
$ws2.Range("C3").Value = 'M20; M35; Ž21A'
$ws2.Range("E3").Value = 2
$ws2.Range("G3").Value = 'Ž21A'
$ws2.Range("H3").Value = 2
$ws2.Range("B4").Value = 'M35'
$ws2.Range("C4").Value = 'M20; M35; Ž21A'
$ws2.Range("E4").Value = 2
$ws2.Range("G4").Value = 'M35'
$ws2.Range("H4").Value = 7
$ws2.Range("B5").Value = 'M20'
$ws2.Range("C5").Value = 'M20; M35; Ž21A'
$ws2.Range("E5").Value = 2
$ws2.Range("G5").Value = 'M20'
$ws2.Range("H5").Value = 9
$ws2.Range("B6").Value = 'M70'
$ws2.Range("C6").Value = 'M65; M70; Ž55'
$ws2.Range("E6").Value = 3
$ws2.Range("G6").Value = 'M70'
$ws2.Range("H6").Value = 3
$ws2.Range("B7").Value = 'M65'
$ws2.Range("C7").Value = 'M65; M70; Ž55'
$ws2.Range("E7").Value = 3
$ws2.Range("G7").Value = 'M65'
$ws2.Range("H7").Value = 4
$ws2.Range("B8").Value = 'Ž55'
$ws2.Range("C8").Value = 'M65; M70; Ž55'
$ws2.Range("E8").Value = 3
$ws2.Range("G8").Value = 'Ž55'
$ws2.Range("H8").Value = 15
$ws2.Range("B9").Value = 'M55'
$ws2.Range("C9").Value = 'M55; Ž45'
$ws2.Range("E9").Value = 5
$ws2.Range("G9").Value = 'M55'
$ws2.Range("H9").Value = 5
$ws2.Range("B10").Value = 'Ž45'
$ws2.Range("C10").Value = 'M55; Ž45'
$ws2.Range("E10").Value = 5
$ws2.Range("G10").Value = 'Ž45'
$ws2.Range("H10").Value = 16
$ws2.Range("B11").Value = 'M45'
$ws2.Range("C11").Value = 'M45; Ž20; Ž35'
$ws2.Range("E11").Value = 6
$ws2.Range("G11").Value = 'M45'
$ws2.Range("H11").Value = 6
$ws2.Range("B12").Value = 'Ž35'
$ws2.Range("C12").Value = 'M45; Ž20; Ž35'
$ws2.Range("E12").Value = 6
$ws2.Range("G12").Value = 'Ž35'
$ws2.Range("H12").Value = 17
$ws2.Range("B13").Value = 'Ž20'
$ws2.Range("C13").Value = 'M45; Ž20; Ž35'
$ws2.Range("E13").Value = 6
$ws2.Range("G13").Value = 'Ž20'
$ws2.Range("H13").Value = 19
$ws2.Range("B14").Value = 'M21B'
$ws2.Range("C14").Value = 'M16; M21B; Open Long'
$ws2.Range("E14").Value = 8
$ws2.Range("G14").Value = 'M21B'
$ws2.Range("H14").Value = 8
$ws2.Range("B15").Value = 'M16'
$ws2.Range("C15").Value = 'M16; M21B; Open Long'
$ws2.Range("E15").Value = 8
$ws2.Range("G15").Value = 'M16'
$ws2.Range("H15").Value = 10
$ws2.Range("B16").Value = 'M16'
$ws2.Range("C16").Value = 'M16; M21B; Open Long'
$ws2.Range("E16").Value = 8
$ws2.Range("G16").Value = 'M16'
$ws2.Range("H16").Value = 10
$ws2.Range("B17").Value = 'OPEN duga'
$ws2.Range("C17").Value = 'M16; M21B; Open Long'
$ws2.Range("E17").Value = 8
$ws2.Range("G17").Value = 'Otvorena duga'
$ws2.Range("H17").Value = 25
$ws2.Range("B18").Value = 'M14'
$ws2.Range("C18").Value = 'M14; Ž14'
$ws2.Range("E18").Value = 11
$ws2.Range("G18").Value = 'M14'
$ws2.Range("H18").Value = 11
$ws2.Range("B19").Value = 'Ž14'
$ws2.Range("C19").Value = 'M14; Ž14'
$ws2.Range("E19").Value = 11
$ws2.Range("G19").Value = 'Ž14'
$ws2.Range("H19").Value = 21
$ws2.Range("B20").Value = 'M12'
$ws2.Range("C20").Value = 'M12; Ž12; Open Short'
$ws2.Range("E20").Value = 12
$ws2.Range("G20").Value = 'M12'
$ws2.Range("H20").Value = 12
$ws2.Range("B21").Value = 'Ž12'
$ws2.Range("C21").Value = 'M12; Ž12; Open Short'
$ws2.Range("E21").Value = 12
$ws2.Range("G21").Value = 'Ž12'
$ws2.Range("H21").Value = 22
$ws2.Range("B22").Value = 'OPEN kratka'
$ws2.Range("C22").Value = 'M12; Ž12; Open Short'
$ws2.Range("E22").Value = 12
$ws2.Range("G22").Value = 'Otvorena kratka'
$ws2.Range("H22").Value = 26
$ws2.Range("B23").Value = 'Ž65'
$ws2.Range("C23").Value = 'Ž65; Ž70'
$ws2.Range("E23").Value = 14
$ws2.Range("G23").Value = 'Ž65'
$ws2.Range("H23").Value = 14
$ws2.Range("B24").Value = 'Ž21B'
$ws2.Range("C24").Value = 'Ž16; Ž21B'
$ws2.Range("E24").Value = 18
$ws2.Range("G24").Value = 'Ž21B'
$ws2.Range("H24").Value = 18
$ws2.Range("B25").Value = 'Ž16'
$ws2.Range("C25").Value = 'Ž16; Ž21B'
$ws2.Range("E25").Value = 18
$ws2.Range("G25").Value = 'Ž16'
$ws2.Range("H25").Value = 20
